$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for columns S (2021) and T (2022), rows 4-14.
# Row 4 = header years; rows 5-14 = data values.
$newData = @{
    4  = @(2021, 2022)
    5  = @(2.5, 2.6)
    6  = @(2.5, 1.8)
    7  = @(1.6, 2.6)
    8  = @(3.6, 1.9)
    9  = @(5.8, 3.9)
    10 = @(1.1000000000000001, 3.2)
    11 = @(1.1000000000000001, 3.3)
    12 = @(5.0999999999999996, 2.5)
    13 = @(2.2999999999999998, 1.9)
    14 = @(2.1, 2.5)
}

foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    # Column S = 19, Column T = 20
    $sCell = $ws.Cells.Item($row, 19)
    $tCell = $ws.Cells.Item($row, 20)

    $sCell.Value = $vals[0]
    $tCell.Value = $vals[1]

    # Match formatting of column R (18) for the same row
    $rCell = $ws.Cells.Item($row, 18)
    $rCell.Copy()
    $sCell.PasteSpecial(-4122) # xlPasteFormats
    $tCell.PasteSpecial(-4122) # xlPasteFormats
}

$excel.CutCopyMode = 0

# Update the selection to match the target state
$ws.Range("V7").Select()
